$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 14.20422986892707
$ws.Cells.Item(2, 3).Value = 10.21172724742363
$ws.Cells.Item(2, 4).Value = 4.832382509627056
$ws.Cells.Item(2, 5).Value = 12.21280148691861
$ws.Cells.Item(2, 6).Value = 24.66609776886348
$ws.Cells.Item(2, 9).Value = 22.11756352582332
$ws.Cells.Item(2, 12).Value = 10.02367379382308
$ws.Cells.Item(2, 13).Value = 14.38219852680219
$ws.Cells.Item(2, 14).Value = 17.83463344637399
$ws.Cells.Item(2, 15).Value = 22.00045793095443

$ws.Cells.Item(3, 2).Value = 13.73990935457991
$ws.Cells.Item(3, 3).Value = 9.980055976118416
$ws.Cells.Item(3, 4).Value = 4.793411487157488
$ws.Cells.Item(3, 5).Value = 12.24915560344527
$ws.Cells.Item(3, 6).Value = 24.64339444581196
$ws.Cells.Item(3, 9).Value = 22.20910413977141
$ws.Cells.Item(3, 12).Value = 10.0313067501568
$ws.Cells.Item(3, 13).Value = 14.29024941157939
$ws.Cells.Item(3, 14).Value = 17.88579679501932
$ws.Cells.Item(3, 15).Value = 22.03916748154745

$ws.Cells.Item(4, 2).Value = 13.44848288234644
$ws.Cells.Item(4, 3).Value = 9.83368689417278
$ws.Cells.Item(4, 4).Value = 4.768982357822719
$ws.Cells.Item(4, 5).Value = 12.27270465952869
$ws.Cells.Item(4, 6).Value = 24.6367315221349
$ws.Cells.Item(4, 9).Value = 22.26998999818573
$ws.Cells.Item(4, 12).Value = 10.0374033365541
$ws.Cells.Item(4, 13).Value = 14.23559042766703
$ws.Cells.Item(4, 14).Value = 17.91897478227963
$ws.Cells.Item(4, 15).Value = 22.06847598540782

$ws.Cells.Item(5, 2).Value = 13.32831459669349
$ws.Cells.Item(5, 3).Value = 9.773055325012376
$ws.Cells.Item(5, 4).Value = 4.75890525872834
$ws.Cells.Item(5, 5).Value = 12.28261052454264
$ws.Cells.Item(5, 6).Value = 24.63584853381371
$ws.Cells.Item(5, 9).Value = 22.29597624348162
$ws.Cells.Item(5, 12).Value = 10.04024278590913
$ws.Cells.Item(5, 13).Value = 14.21378646875386
$ws.Cells.Item(5, 14).Value = 17.93293948704045
$ws.Cells.Item(5, 15).Value = 22.08180949640498

$ws.Cells.Item(6, 2).Value = 13.30828123642048
$ws.Cells.Item(6, 3).Value = 9.762929575081362
$ws.Cells.Item(6, 4).Value = 4.757224708088927
$ws.Cells.Item(6, 5).Value = 12.28427409769026
$ws.Cells.Item(6, 6).Value = 24.63581261118858
$ws.Cells.Item(6, 9).Value = 22.30036213564803
$ws.Cells.Item(6, 12).Value = 10.04073573121935
$ws.Cells.Item(6, 13).Value = 14.21019483447679
$ws.Cells.Item(6, 14).Value = 17.93528518519691
$ws.Cells.Item(6, 15).Value = 22.08410736927178

$ws.Cells.Item(7, 2).Value = 13.44686770607912
$ws.Cells.Item(7, 3).Value = 9.832873113424817
$ws.Cells.Item(7, 4).Value = 4.76884694348175
$ws.Cells.Item(7, 5).Value = 12.27283699957019
$ws.Cells.Item(7, 6).Value = 24.63671219345723
$ws.Cells.Item(7, 9).Value = 22.27033570379059
$ws.Cells.Item(7, 12).Value = 10.03744019219969
$ws.Cells.Item(7, 13).Value = 14.23529444567842
$ws.Cells.Item(7, 14).Value = 17.91916131450795
$ws.Cells.Item(7, 15).Value = 22.06865018209569

$ws.Cells.Item(8, 2).Value = 14.0455511585731
$ws.Cells.Item(8, 3).Value = 10.13273245453537
$ws.Cells.Item(8, 4).Value = 4.819051554051577
$ws.Cells.Item(8, 5).Value = 12.22508208881229
$ws.Cells.Item(8, 6).Value = 24.65676121854674
$ws.Cells.Item(8, 9).Value = 22.14815432106385
$ws.Cells.Item(8, 12).Value = 10.02601333122909
$ws.Cells.Item(8, 13).Value = 14.35013218531034
$ws.Cells.Item(8, 14).Value = 17.85190913603167
$ws.Cells.Item(8, 15).Value = 22.01265323368988

$ws.Cells.Item(9, 2).Value = 15.16207796500663
$ws.Cells.Item(9, 3).Value = 10.68594802291027
$ws.Cells.Item(9, 4).Value = 4.913361899823452
$ws.Cells.Item(9, 5).Value = 12.14113898361764
$ws.Cells.Item(9, 6).Value = 24.7536407858855
$ws.Cells.Item(9, 9).Value = 21.94577388288311
$ws.Cells.Item(9, 12).Value = 10.01476814830278
$ws.Cells.Item(9, 13).Value = 14.58878220473857
$ws.Cells.Item(9, 14).Value = 17.73397659386125
$ws.Cells.Item(9, 15).Value = 21.9469297556047

$ws.Cells.Item(10, 2).Value = 15.93870137602916
$ws.Cells.Item(10, 3).Value = 11.06859139950381
$ws.Cells.Item(10, 4).Value = 4.979908562591143
$ws.Cells.Item(10, 5).Value = 12.08533221235153
$ws.Cells.Item(10, 6).Value = 24.85956691224339
$ws.Cells.Item(10, 9).Value = 21.81989157081367
$ws.Cells.Item(10, 12).Value = 10.01327586716003
$ws.Cells.Item(10, 13).Value = 14.77119305268255
$ws.Cells.Item(10, 14).Value = 17.65577435098562
$ws.Cells.Item(10, 15).Value = 21.92565899360954

$ws.Cells.Item(11, 2).Value = 16.28098764165281
$ws.Cells.Item(11, 3).Value = 11.23701457216906
$ws.Cells.Item(11, 4).Value = 5.009541439059879
$ws.Cells.Item(11, 5).Value = 12.06120753507879
$ws.Cells.Item(11, 6).Value = 24.91519212377967
$ws.Cells.Item(11, 9).Value = 21.76760207721424
$ws.Cells.Item(11, 12).Value = 10.01405776597997
$ws.Cells.Item(11, 13).Value = 14.85546346271508
$ws.Cells.Item(11, 14).Value = 17.62201874855022
$ws.Cells.Item(11, 15).Value = 21.921867684301

$ws.Cells.Item(12, 2).Value = 16.40891090955349
$ws.Cells.Item(12, 3).Value = 11.29994440079854
$ws.Cells.Item(12, 4).Value = 5.020666903480038
$ws.Cells.Item(12, 5).Value = 12.05225285282173
$ws.Cells.Item(12, 6).Value = 24.93731375038714
$ws.Cells.Item(12, 9).Value = 21.74851896836594
$ws.Cells.Item(12, 12).Value = 10.01456296636933
$ws.Cells.Item(12, 13).Value = 14.88753860184055
$ws.Cells.Item(12, 14).Value = 17.60949703548171
$ws.Cells.Item(12, 15).Value = 21.92127880817942

$ws.Cells.Item(13, 2).Value = 16.38143727616387
$ws.Cells.Item(13, 3).Value = 11.28642962240463
$ws.Cells.Item(13, 4).Value = 5.018275160398827
$ws.Cells.Item(13, 5).Value = 12.05417337305327
$ws.Cells.Item(13, 6).Value = 24.93250263758991
$ws.Cells.Item(13, 9).Value = 21.75259689559649
$ws.Cells.Item(13, 12).Value = 10.01444487801909
$ws.Cells.Item(13, 13).Value = 14.88062371322844
$ws.Cells.Item(13, 14).Value = 17.61218222323991
$ws.Cells.Item(13, 15).Value = 21.92136796829704

$ws.Cells.Item(14, 2).Value = 16.29154648791466
$ws.Cells.Item(14, 3).Value = 11.24220905000536
$ws.Cells.Item(14, 4).Value = 5.010458680992097
$ws.Cells.Item(14, 5).Value = 12.06046720815289
$ws.Cells.Item(14, 6).Value = 24.91699096351458
$ws.Cells.Item(14, 9).Value = 21.76601769919435
$ws.Cells.Item(14, 12).Value = 10.0140951438999
$ws.Cells.Item(14, 13).Value = 14.85809914748871
$ws.Cells.Item(14, 14).Value = 17.62098335716274
$ws.Cells.Item(14, 15).Value = 21.92180226460585

$ws.Cells.Item(15, 2).Value = 16.23626222024955
$ws.Cells.Item(15, 3).Value = 11.21501114746606
$ws.Cells.Item(15, 4).Value = 5.005658263153424
$ws.Cells.Item(15, 5).Value = 12.06434589058189
$ws.Cells.Item(15, 6).Value = 24.90762695309823
$ws.Cells.Item(15, 9).Value = 21.77433187672968
$ws.Cells.Item(15, 12).Value = 10.01390812482154
$ws.Cells.Item(15, 13).Value = 14.84432287560516
$ws.Cells.Item(15, 14).Value = 17.62640824737828
$ws.Cells.Item(15, 15).Value = 21.92217857017674

$ws.Cells.Item(16, 2).Value = 15.91610097727517
$ws.Cells.Item(16, 3).Value = 11.0574677695656
$ws.Cells.Item(16, 4).Value = 4.977958780305001
$ws.Cells.Item(16, 5).Value = 12.08693415522646
$ws.Cells.Item(16, 6).Value = 24.85608036381556
$ws.Cells.Item(16, 9).Value = 21.82340912635992
$ws.Cells.Item(16, 12).Value = 10.01325408370963
$ws.Cells.Item(16, 13).Value = 14.7657098942164
$ws.Cells.Item(16, 14).Value = 17.65801689563104
$ws.Cells.Item(16, 15).Value = 21.92602524581245

$ws.Cells.Item(17, 2).Value = 15.71679118586167
$ws.Cells.Item(17, 3).Value = 10.95934741867641
$ws.Cells.Item(17, 4).Value = 4.9607994461925
$ws.Cells.Item(17, 5).Value = 12.1011140928528
$ws.Cells.Item(17, 6).Value = 24.82635537334022
$ws.Cells.Item(17, 9).Value = 21.85479234938812
$ws.Cells.Item(17, 12).Value = 10.0132263406912
$ws.Cells.Item(17, 13).Value = 14.7177985035184
$ws.Cells.Item(17, 14).Value = 17.67787311179012
$ws.Cells.Item(17, 15).Value = 21.92989290939618

$ws.Cells.Item(18, 2).Value = 15.60112352790553
$ws.Cells.Item(18, 3).Value = 10.90238202693039
$ws.Cells.Item(18, 4).Value = 4.950869883469499
$ws.Cells.Item(18, 5).Value = 12.10938885315845
$ws.Cells.Item(18, 6).Value = 24.80995936166945
$ws.Cells.Item(18, 9).Value = 21.87331127762756
$ws.Cells.Item(18, 12).Value = 10.01334793770691
$ws.Cells.Item(18, 13).Value = 14.69036395884098
$ws.Cells.Item(18, 14).Value = 17.68946513865094
$ws.Cells.Item(18, 15).Value = 21.93267140182965

$ws.Cells.Item(19, 2).Value = 15.56178714121106
$ws.Cells.Item(19, 3).Value = 10.88300481474309
$ws.Cells.Item(19, 4).Value = 4.947497724728716
$ws.Cells.Item(19, 5).Value = 12.11221097571167
$ws.Cells.Item(19, 6).Value = 24.80452869660922
$ws.Cells.Item(19, 9).Value = 21.87966179710945
$ws.Cells.Item(19, 12).Value = 10.0134127606522
$ws.Cells.Item(19, 13).Value = 14.68109685017613
$ws.Cells.Item(19, 14).Value = 17.69341943697648
$ws.Cells.Item(19, 15).Value = 21.93370725979189

$ws.Cells.Item(20, 2).Value = 15.73811549132942
$ws.Cells.Item(20, 3).Value = 10.96984754069166
$ws.Cells.Item(20, 4).Value = 4.962632321507059
$ws.Cells.Item(20, 5).Value = 12.09959232023515
$ws.Cells.Item(20, 6).Value = 24.82944717934257
$ws.Cells.Item(20, 9).Value = 21.85140308483655
$ws.Cells.Item(20, 12).Value = 10.01321506419422
$ws.Cells.Item(20, 13).Value = 14.72288619906719
$ws.Cells.Item(20, 14).Value = 17.67574166489586
$ws.Cells.Item(20, 15).Value = 21.92942385585267

$ws.Cells.Item(21, 2).Value = 16.3179963620722
$ws.Cells.Item(21, 3).Value = 11.25522100729369
$ws.Cells.Item(21, 4).Value = 5.012757202827734
$ws.Cells.Item(21, 5).Value = 12.05861365447852
$ws.Cells.Item(21, 6).Value = 24.92151852108579
$ws.Cells.Item(21, 9).Value = 21.76205618432766
$ws.Cells.Item(21, 12).Value = 10.01419220197534
$ws.Cells.Item(21, 13).Value = 14.86471088904216
$ws.Cells.Item(21, 14).Value = 17.61839117873711
$ws.Cells.Item(21, 15).Value = 21.92165171745262

$ws.Cells.Item(22, 2).Value = 16.68707065303462
$ws.Cells.Item(22, 3).Value = 11.43676914436546
$ws.Cells.Item(22, 4).Value = 5.044956198173855
$ws.Cells.Item(22, 5).Value = 12.03288531370434
$ws.Cells.Item(22, 6).Value = 24.98784998297331
$ws.Cells.Item(22, 9).Value = 21.70784773429583
$ws.Cells.Item(22, 12).Value = 10.01604918594008
$ws.Cells.Item(22, 13).Value = 14.95834581128244
$ws.Cells.Item(22, 14).Value = 17.58242912699108
$ws.Cells.Item(22, 15).Value = 21.92150787673156

$ws.Cells.Item(23, 2).Value = 16.49102901586571
$ws.Cells.Item(23, 3).Value = 11.34033871521429
$ws.Cells.Item(23, 4).Value = 5.027823542241048
$ws.Cells.Item(23, 5).Value = 12.0465208325961
$ws.Cells.Item(23, 6).Value = 24.9518886550227
$ws.Cells.Item(23, 9).Value = 21.73639609841329
$ws.Cells.Item(23, 12).Value = 10.01494693143521
$ws.Cells.Item(23, 13).Value = 14.90829195270242
$ws.Cells.Item(23, 14).Value = 17.60148393666336
$ws.Cells.Item(23, 15).Value = 21.92113299575967

$ws.Cells.Item(24, 2).Value = 15.7284781338427
$ws.Cells.Item(24, 3).Value = 10.96510216049133
$ws.Cells.Item(24, 4).Value = 4.961803878971245
$ws.Cells.Item(24, 5).Value = 12.10027993193433
$ws.Cells.Item(24, 6).Value = 24.82804721341592
$ws.Cells.Item(24, 9).Value = 21.85293388787295
$ws.Cells.Item(24, 12).Value = 10.0132197337858
$ws.Cells.Item(24, 13).Value = 14.7205857065034
$ws.Cells.Item(24, 14).Value = 17.67670474241859
$ws.Cells.Item(24, 15).Value = 21.92963418645103

$ws.Cells.Item(25, 2).Value = 14.86713266287922
$ws.Cells.Item(25, 3).Value = 10.54029433400176
$ws.Cells.Item(25, 4).Value = 4.888314543836252
$ws.Cells.Item(25, 5).Value = 12.16281408709203
$ws.Cells.Item(25, 6).Value = 24.72129841209278
$ws.Cells.Item(25, 9).Value = 21.99652661952814
$ws.Cells.Item(25, 12).Value = 10.01661855879601
$ws.Cells.Item(25, 13).Value = 14.5228951358727
$ws.Cells.Item(25, 14).Value = 17.76439351595178
$ws.Cells.Item(25, 15).Value = 21.95997336861909
